$d = $word.ActiveDocument

# Locate the paragraph that holds the secretary's name line
# ("Секретарь инвестиционного комитета <TAB> Оспанова I.Е.").
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*Оспанова*") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not find the paragraph containing 'Оспанова'"
}

# Replace the whole paragraph with the corrected OOXML:
#  - paragraph-mark language switched from ru-RU to en-US
#  - the "Оспанова I.Е." name (plus its spell-check proof markers)
#    collapsed into a single "ICADMIN" run, keeping the same bold /
#    size / en-US run formatting the trailing initials already had
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="5D30D1DE" w14:textId="2985D66E" w:rsidR="001C6721" w:rsidRDefault="001C6721" w:rsidP="001C6721"><w:pPr><w:pStyle w:val="1"/><w:shd w:val="clear" w:color="auto" w:fill="auto"/><w:tabs><w:tab w:val="left" w:pos="7230"/><w:tab w:val="left" w:leader="underscore" w:pos="8858"/></w:tabs><w:spacing w:line="240" w:lineRule="auto"/><w:ind w:firstLine="403"/><w:jc w:val="both"/><w:rPr><w:bCs/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r w:rsidRPr="00483757"><w:rPr><w:b/><w:bCs/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="ru-RU"/></w:rPr><w:t>Секр</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="ru-RU"/></w:rPr><w:t>етарь инвестиционного комитета</w:t></w:r><w:r w:rsidRPr="006B78D7"><w:rPr><w:bCs/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="ru-RU"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:bCs/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="ru-RU"/></w:rPr><w:tab/></w:r><w:r w:rsidR="007010DF"><w:rPr><w:b/><w:bCs/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="en-US"/></w:rPr><w:t>ICADMIN</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'

$target.Range.InsertXML($xml)
